$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on cells whose new numeric-looking values must remain text strings
$textCells = @("D5","D6","D7","D9","D10","D12","D13","D15","D17","D19","D20","D21","D23","D24","D25","D26","D27","D29","D30","D31","D32","D33","D35","D36","D38","D39","D40","D41","D44","D45","D49","D51")
foreach ($c in $textCells) {
    $ws.Range($c).NumberFormat = "@"
}

$ws.Range("D2").Value = "50.898.02"
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").Value = "2.900.70"
$ws.Range("E3").Value = "  -0.29%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "368.35"
$ws.Range("E5").Value = "  +5.85%  "
$ws.Range("D6").Value = "103.14"
$ws.Range("E6").Value = "  -2.07%  "
$ws.Range("D7").Value = "0.539"
$ws.Range("E7").Value = "  -2.10%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "0.583"
$ws.Range("E9").Value = "  -3.32%  "
$ws.Range("D10").Value = "36.37"
$ws.Range("E10").Value = "  -2.58%  "
$ws.Range("E11").Value = "  +0.85%  "
$ws.Range("D12").Value = "0.0830"
$ws.Range("E12").Value = "  -1.11%  "
$ws.Range("D13").Value = "18.21"
$ws.Range("E13").Value = "  -2.69%  "
$ws.Range("D14").Value = "3.356.74"
$ws.Range("E14").Value = "  -0.28%  "
$ws.Range("D15").Value = "7.36"
$ws.Range("E15").Value = "  -1.81%  "
$ws.Range("D16").Value = "2.902.17"
$ws.Range("E16").Value = "  -0.02%  "
$ws.Range("D17").Value = "0.923"
$ws.Range("E17").Value = "  -2.81%  "
$ws.Range("D18").Value = "50.887.42"
$ws.Range("E18").Value = "  -0.45%  "
$ws.Range("D19").Value = "3.21"
$ws.Range("E19").Value = "  -5.62%  "
$ws.Range("D20").Value = "7.16"
$ws.Range("E20").Value = "  -2.25%  "
$ws.Range("D21").Value = "12.82"
$ws.Range("E21").Value = "  -3.51%  "
$ws.Range("D22").Value = "0.0₃0938"
$ws.Range("E22").Value = "  -1.91%  "
$ws.Range("D23").Value = "68.24"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").Value = "258.38"
$ws.Range("E24").Value = "  -0.71%  "
$ws.Range("D25").Value = "2.67"
$ws.Range("E25").Value = "  -0.31%  "
$ws.Range("D26").Value = "4.32"
$ws.Range("E26").Value = "  +3.38%  "
$ws.Range("D27").Value = "0.174"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E28").Value = "  -0.07%  "
$ws.Range("D29").Value = "25.57"
$ws.Range("E29").Value = "  -2.34%  "
$ws.Range("D30").Value = "6.94"
$ws.Range("E30").Value = "  -5.68%  "
$ws.Range("D31").Value = "0.101"
$ws.Range("E31").Value = "  -1.52%  "
$ws.Range("D32").Value = "6.17"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("D33").Value = "9.84"
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("E34").Value = "  -0.30%  "
$ws.Range("D35").Value = "34.35"
$ws.Range("E35").Value = "  -2.64%  "
$ws.Range("D36").Value = "50.83"
$ws.Range("E37").Value = "  +0.57%  "
$ws.Range("D38").Value = "0.0420"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").Value = "2.99"
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("D40").Value = "2.62"
$ws.Range("E40").Value = "  +0.50%  "
$ws.Range("D41").Value = "16.95"
$ws.Range("E41").Value = "  -3.01%  "
$ws.Range("E42").Value = "  -4.68%  "
$ws.Range("E43").Value = "  -2.47%  "
$ws.Range("D44").Value = "21.99"
$ws.Range("E44").Value = "  -1.25%  "
$ws.Range("D45").Value = "118.11"
$ws.Range("E45").Value = "  -1.32%  "
$ws.Range("E46").Value = "  -2.25%  "
$ws.Range("D47").Value = "2.010.24"
$ws.Range("E47").Value = "  -3.83%  "
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("D49").Value = "3.14"
$ws.Range("E49").Value = "  -4.02%  "
$ws.Range("D50").Value = "3.200.74"
$ws.Range("E50").Value = "  +0.23%  "
$ws.Range("D51").Value = "0.237"
$ws.Range("E51").Value = "  +0.01%  "
